# Update country list & provincias Spain (covid-19 "Pais" dashboard)
# - refresh "data as of" timestamp
# - fix several shared-string country-name cells that had drifted out of
#   their sorted position (values themselves are untouched, only which
#   country label sits in a given row)
# - refresh the per-country case/recovered/death counters with the newer pull

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- "Datos actualizados a ..." banner -----------------------------------
$ws.Range("A1").Value = 'Datos actualizados a 23 de Marzo de 2020 a las 12:46'

# --- Country-name cells that need to point at the right label ------------

$ws.Range("A59").Value = 'Irak'
$ws.Range("A60").Value = 'Libano'
$ws.Range("A61").Value = 'Colombia'

$ws.Range("A76").Value = 'Republica de Macedonia'
$ws.Range("A77").Value = 'Costa Rica'
$ws.Range("A78").Value = 'Bosnia y Herzegovina'
$ws.Range("A79").Value = 'Marruecos'
$ws.Range("A80").Value = 'Vietnam'
$ws.Range("A81").Value = 'Islas Feroe'

$ws.Range("A84").Value = 'Malta'
$ws.Range("A85").Value = 'Nueva Zelanda'
$ws.Range("A86").Value = 'Burkina Faso'
$ws.Range("A87").Value = 'Republica de Chipre'
$ws.Range("A88").Value = 'Moldavia'
$ws.Range("A89").Value = 'Brunei'
$ws.Range("A90").Value = 'Albania'
$ws.Range("A91").Value = 'Sri Lanka'
$ws.Range("A92").Value = 'Camboya'
$ws.Range("A93").Value = 'Bielorrusia'
$ws.Range("A94").Value = 'Venezuela'
$ws.Range("A95").Value = 'Tunez'

$ws.Range("A114").Value = 'Puerto Rico'
$ws.Range("A115").Value = 'Consejo Danes para los Refugiados'
$ws.Range("A116").Value = 'Guam'
$ws.Range("A117").Value = 'Mauricio'

$ws.Range("A119").Value = 'Bolivia'
$ws.Range("A120").Value = 'Costa de Marfil'
$ws.Range("A121").Value = 'Ghana'
$ws.Range("A122").Value = 'Macao'
$ws.Range("A123").Value = 'Monaco'
$ws.Range("A124").Value = 'Montenegro'
$ws.Range("A125").Value = 'Paraguay'

$ws.Range("A132").Value = 'Togo'
$ws.Range("A133").Value = 'Polinesia Francesa'
$ws.Range("A134").Value = 'Barbados'

# --- Refreshed statistics (Casos totales / Nuevos casos / Casos activos /
#     Recuperados / Casos criticos / Muertes hoy / Muertes) ----------------

# Row 10
$ws.Range("B10").Value = 16481
$ws.Range("C10").Value = 463
$ws.Range("E10").Value = 13607

# Row 12
$ws.Range("E12").Value = 7996
$ws.Range("G12").Value = 9
$ws.Range("H12").Value = 107

# Row 15
$ws.Range("B15").Value = 3784
$ws.Range("C15").Value = 202
$ws.Range("E15").Value = 3759

# Row 20
$ws.Range("B20").Value = 1619
$ws.Range("C20").Value = 73
$ws.Range("E20").Value = 1592

# Row 36
$ws.Range("B36").Value = 684
$ws.Range("C36").Value = 50
$ws.Range("E36").Value = 663
$ws.Range("G36").Value = 1
$ws.Range("H36").Value = 8

# Row 38
$ws.Range("E38").Value = 588
$ws.Range("G38").Value = 2
$ws.Range("H38").Value = 17

# Row 40
$ws.Range("F40").Value = 15

# Row 46
$ws.Range("D46").Value = 17
$ws.Range("E46").Value = 420

# Row 59
$ws.Range("B59").Value = 266
$ws.Range("C59").Value = 33
$ws.Range("D59").Value = 62
$ws.Range("E59").Value = 181
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 3
$ws.Range("H59").Value = 23

# Row 60
$ws.Range("B60").Value = 256
$ws.Range("C60").Value = 8
$ws.Range("D60").Value = 8
$ws.Range("E60").Value = 244
$ws.Range("F60").Value = 4
$ws.Range("H60").Value = 4

# Row 61
$ws.Range("B61").Value = 235
$ws.Range("C61").Value = 4
$ws.Range("D61").Value = 3
$ws.Range("E61").Value = 230
$ws.Range("H61").Value = 2

# Row 76
$ws.Range("B76").Value = 136
$ws.Range("C76").Value = 21
$ws.Range("D76").Value = 1
$ws.Range("E76").Value = 133
$ws.Range("F76").Value = 1
$ws.Range("G76").Value = 1

# Row 77
$ws.Range("B77").Value = 134
$ws.Range("C77").Value = 0
$ws.Range("E77").Value = 130
$ws.Range("F77").Value = 2
$ws.Range("H77").Value = 2

# Row 78
$ws.Range("B78").Value = 128
$ws.Range("C78").Value = 2
$ws.Range("D78").Value = 2
$ws.Range("E78").Value = 125
$ws.Range("H78").Value = 1

# Row 79
$ws.Range("B79").Value = 122
$ws.Range("C79").Value = 7
$ws.Range("D79").Value = 3
$ws.Range("E79").Value = 115
$ws.Range("F79").Value = 1
$ws.Range("H79").Value = 4

# Row 80
$ws.Range("B80").Value = 121
$ws.Range("C80").Value = 8
$ws.Range("D80").Value = 17
$ws.Range("F80").Value = 2

# Row 81
$ws.Range("B81").Value = 118
$ws.Range("C81").Value = 3
$ws.Range("D81").Value = 14
$ws.Range("E81").Value = 104
$ws.Range("F81").Value = 0
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 0

# Row 84
$ws.Range("B84").Value = 107
$ws.Range("C84").Value = 17
$ws.Range("D84").Value = 2
$ws.Range("E84").Value = 105
$ws.Range("F84").Value = 1

# Row 85
$ws.Range("B85").Value = 102
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 102
$ws.Range("F85").Value = 0
$ws.Range("H85").Value = 0

# Row 86
$ws.Range("B86").Value = 99
$ws.Range("C86").Value = 24
$ws.Range("D86").Value = 5
$ws.Range("E86").Value = 90
$ws.Range("F86").Value = 0
$ws.Range("H86").Value = 4

# Row 87
$ws.Range("B87").Value = 95
$ws.Range("C87").Value = 0
$ws.Range("D87").Value = 3
$ws.Range("E87").Value = 91
$ws.Range("F87").Value = 3
$ws.Range("H87").Value = 1

# Row 88
$ws.Range("B88").Value = 94
$ws.Range("E88").Value = 91
$ws.Range("F88").Value = 3
$ws.Range("H88").Value = 1

# Row 89
$ws.Range("B89").Value = 91
$ws.Range("C89").Value = 3
$ws.Range("E89").Value = 89
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 0

# Row 90
$ws.Range("B90").Value = 89
$ws.Range("C90").Value = 0
$ws.Range("D90").Value = 2
$ws.Range("E90").Value = 82
$ws.Range("G90").Value = 3
$ws.Range("H90").Value = 5

# Row 91
$ws.Range("B91").Value = 87
$ws.Range("C91").Value = 5
$ws.Range("D91").Value = 3
$ws.Range("F91").Value = 2

# Row 92
$ws.Range("B92").Value = 86
$ws.Range("C92").Value = 2
$ws.Range("D92").Value = 2
$ws.Range("E92").Value = 84
$ws.Range("F92").Value = 1

# Row 93
$ws.Range("B93").Value = 81
$ws.Range("C93").Value = 5
$ws.Range("D93").Value = 22
$ws.Range("E93").Value = 59
$ws.Range("F93").Value = 0

# Row 94
$ws.Range("B94").Value = 77
$ws.Range("C94").Value = 7
$ws.Range("D94").Value = 15
$ws.Range("E94").Value = 62
$ws.Range("F94").Value = 2
$ws.Range("H94").Value = 0

# Row 95
$ws.Range("D95").Value = 1
$ws.Range("E95").Value = 71
$ws.Range("F95").Value = 7
$ws.Range("H95").Value = 3

# Row 113
$ws.Range("D113").Value = 5
$ws.Range("E113").Value = 25

# Row 114
$ws.Range("B114").Value = 31
$ws.Range("C114").Value = 8
$ws.Range("G114").Value = 1
$ws.Range("H114").Value = 2

# Row 115
$ws.Range("B115").Value = 30
$ws.Range("C115").Value = 0
$ws.Range("E115").Value = 29

# Row 116
$ws.Range("B116").Value = 29
$ws.Range("C116").Value = 2
$ws.Range("E116").Value = 28
$ws.Range("F116").Value = 0
$ws.Range("H116").Value = 1

# Row 117
$ws.Range("B117").Value = 28
$ws.Range("C117").Value = 0
$ws.Range("E117").Value = 26
$ws.Range("F117").Value = 1
$ws.Range("H117").Value = 2

# Row 118
$ws.Range("B118").Value = 27
$ws.Range("C118").Value = 1
$ws.Range("E118").Value = 27

# Row 119
$ws.Range("B119").Value = 27
$ws.Range("C119").Value = 3
$ws.Range("D119").Value = 0
$ws.Range("E119").Value = 27

# Row 120
$ws.Range("B120").Value = 25
$ws.Range("C120").Value = 11
$ws.Range("D120").Value = 2
$ws.Range("H120").Value = 0

# Row 121
$ws.Range("C121").Value = 1
$ws.Range("D121").Value = 0
$ws.Range("E121").Value = 23
$ws.Range("H121").Value = 1

# Row 122
$ws.Range("B122").Value = 24
$ws.Range("C122").Value = 2
$ws.Range("D122").Value = 10
$ws.Range("E122").Value = 14

# Row 123
$ws.Range("D123").Value = 1
$ws.Range("E123").Value = 22
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 0

# Row 124
$ws.Range("C124").Value = 1
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 1

# Row 125
$ws.Range("C125").Value = 0
$ws.Range("F125").Value = 1
$ws.Range("G125").Value = 0

# Row 132
$ws.Range("C132").Value = 2

# Row 133
$ws.Range("B133").Value = 18
$ws.Range("C133").Value = 0
$ws.Range("E133").Value = 18

# Row 134
$ws.Range("B134").Value = 17
$ws.Range("C134").Value = 3
$ws.Range("E134").Value = 17
